$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Btc"
$ws.Range("C2").Value = "Egfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.229324
$ws.Range("H2").Value = 0.687972
$ws.Range("I2").Value = 0.1345276800989288
$ws.Range("J2").Value = 0.1345276800989288
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.307106666666667
$ws.Range("N2").Value = 3.92132
$ws.Range("O2").Value = 0.01256263154946851
$ws.Range("P2").Value = 0.01256263154946851
$ws.Range("Q2").Value = 0.2997509292266667
$ws.Range("R2").Value = 2.69775836304
$ws.Range("S2").Value = 0.00169002167828761
$ws.Range("T2").Value = 0.00169002167828761

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Btc"
$ws.Range("C3").Value = "Egfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.229324
$ws.Range("H3").Value = 0.687972
$ws.Range("I3").Value = 0.1345276800989288
$ws.Range("J3").Value = 0.1345276800989288
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.77105616682495
$ws.Range("P3").Value = 0.77105616682495
$ws.Range("Q3").Value = 18.397801573788
$ws.Range("R3").Value = 165.580214164092
$ws.Range("S3").Value = 0.1037283973489331
$ws.Range("T3").Value = 0.1037283973489331

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Btc"
$ws.Range("C4").Value = "Egfr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.229324
$ws.Range("H4").Value = 0.687972
$ws.Range("I4").Value = 0.1345276800989288
$ws.Range("J4").Value = 0.1345276800989288
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.51385866666667
$ws.Range("N4").Value = 67.541576
$ws.Range("O4").Value = 0.2163812016255815
$ws.Range("P4").Value = 0.2163812016255815
$ws.Range("Q4").Value = 5.162968124874667
$ws.Range("R4").Value = 46.466713123872
$ws.Range("S4").Value = 0.02910926107170803
$ws.Range("T4").Value = 0.02910926107170804

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Btc"
$ws.Range("C5").Value = "Egfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.475336333333334
$ws.Range("H5").Value = 4.426009000000001
$ws.Range("I5").Value = 0.8654723199010712
$ws.Range("J5").Value = 0.8654723199010712
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.307106666666667
$ws.Range("N5").Value = 3.92132
$ws.Range("O5").Value = 0.01256263154946851
$ws.Range("P5").Value = 0.01256263154946851
$ws.Range("Q5").Value = 1.928421956875556
$ws.Range("R5").Value = 17.35579761188
$ws.Range("S5").Value = 0.0108726098711809
$ws.Range("T5").Value = 0.0108726098711809

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Btc"
$ws.Range("C6").Value = "Egfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.475336333333334
$ws.Range("H6").Value = 4.426009000000001
$ws.Range("I6").Value = 0.8654723199010712
$ws.Range("J6").Value = 0.8654723199010712
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.77105616682495
$ws.Range("P6").Value = 0.77105616682495
$ws.Range("Q6").Value = 118.360682332711
$ws.Range("R6").Value = 1065.246140994399
$ws.Range("S6").Value = 0.6673277694760169
$ws.Range("T6").Value = 0.6673277694760169

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Btc"
$ws.Range("C7").Value = "Egfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.475336333333334
$ws.Range("H7").Value = 4.426009000000001
$ws.Range("I7").Value = 0.8654723199010712
$ws.Range("J7").Value = 0.8654723199010712
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.51385866666667
$ws.Range("N7").Value = 67.541576
$ws.Range("O7").Value = 0.2163812016255815
$ws.Range("P7").Value = 0.2163812016255815
$ws.Range("Q7").Value = 33.2155136944649
$ws.Range("R7").Value = 298.9396232501841
$ws.Range("S7").Value = 0.1872719405538734
$ws.Range("T7").Value = 0.1872719405538734
